$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(324).Delete()

Write-Host ("A324=" + $ws.Range("A324").Value())
Write-Host ("B324=" + $ws.Range("B324").Value())
Write-Host ("C324=" + $ws.Range("C324").Value())
Write-Host ("D324=" + $ws.Range("D324").Value())
Write-Host ("E324=" + $ws.Range("E324").Value())
Write-Host ("Dimension=" + $ws.UsedRange.Address())
